$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.023873
$ws.Range("H2").Value = 0.071619
$ws.Range("I2").Value = 0.02747901635872243
$ws.Range("J2").Value = 0.02747901635872243
$ws.Range("M2").Value = 15.47987166666667
$ws.Range("N2").Value = 46.439615
$ws.Range("O2").Value = 0.7960757698994193
$ws.Range("P2").Value = 0.7960757698994194
$ws.Range("Q2").Value = 0.3695509762983334
$ws.Range("R2").Value = 3.325958786685
$ws.Range("S2").Value = 0.02187537910384869
$ws.Range("T2").Value = 0.02187537910384869
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.023873
$ws.Range("H3").Value = 0.071619
$ws.Range("I3").Value = 0.02747901635872243
$ws.Range("J3").Value = 0.02747901635872243
$ws.Range("O3").Value = 0.05100527512565552
$ws.Range("P3").Value = 0.05100527512565553
$ws.Range("Q3").Value = 0.023677456257
$ws.Range("R3").Value = 0.213097106313
$ws.Range("S3").Value = 0.001401574789559026
$ws.Range("T3").Value = 0.001401574789559026
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.023873
$ws.Range("H4").Value = 0.071619
$ws.Range("I4").Value = 0.02747901635872243
$ws.Range("J4").Value = 0.02747901635872243
$ws.Range("M4").Value = 0.171678
$ws.Range("N4").Value = 0.515034
$ws.Range("O4").Value = 0.008828800326496623
$ws.Range("P4").Value = 0.008828800326496624
$ws.Range("Q4").Value = 0.004098468894
$ws.Range("R4").Value = 0.036886220046
$ws.Range("S4").Value = 0.0002426067485996946
$ws.Range("T4").Value = 0.0002426067485996946
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.023873
$ws.Range("H5").Value = 0.071619
$ws.Range("I5").Value = 0.02747901635872243
$ws.Range("J5").Value = 0.02747901635872243
$ws.Range("M5").Value = 2.688466333333334
$ws.Range("N5").Value = 8.065399000000001
$ws.Range("O5").Value = 0.1382584398787761
$ws.Range("P5").Value = 0.1382584398787761
$ws.Range("Q5").Value = 0.06418175677566669
$ws.Range("R5").Value = 0.5776358109810001
$ws.Range("S5").Value = 0.003799205931160328
$ws.Range("T5").Value = 0.003799205931160328
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.023873
$ws.Range("H6").Value = 0.071619
$ws.Range("I6").Value = 0.02747901635872243
$ws.Range("J6").Value = 0.02747901635872243
$ws.Range("M6").Value = 0.113399
$ws.Range("N6").Value = 0.340197
$ws.Range("O6").Value = 0.005831714769652435
$ws.Range("P6").Value = 0.005831714769652436
$ws.Range("Q6").Value = 0.002707174327000001
$ws.Range("R6").Value = 0.024364568943
$ws.Range("S6").Value = 0.0001602497855546824
$ws.Range("T6").Value = 0.0001602497855546825
$ws.Range("I7").Value = 0.04107483513127341
$ws.Range("J7").Value = 0.04107483513127341
$ws.Range("M7").Value = 15.47987166666667
$ws.Range("N7").Value = 46.439615
$ws.Range("O7").Value = 0.7960757698994193
$ws.Range("P7").Value = 0.7960757698994194
$ws.Range("Q7").Value = 0.5523940604677777
$ws.Range("R7").Value = 4.971546544210001
$ws.Range("S7").Value = 0.03269868100062019
$ws.Range("T7").Value = 0.03269868100062021
$ws.Range("I8").Value = 0.04107483513127341
$ws.Range("J8").Value = 0.04107483513127341
$ws.Range("O8").Value = 0.05100527512565552
$ws.Range("P8").Value = 0.05100527512565553
$ws.Range("S8").Value = 0.002095033266611541
$ws.Range("T8").Value = 0.002095033266611542
$ws.Range("I9").Value = 0.04107483513127341
$ws.Range("J9").Value = 0.04107483513127341
$ws.Range("M9").Value = 0.171678
$ws.Range("N9").Value = 0.515034
$ws.Range("O9").Value = 0.008828800326496623
$ws.Range("P9").Value = 0.008828800326496624
$ws.Range("Q9").Value = 0.006126272203999999
$ws.Range("R9").Value = 0.055136449836
$ws.Range("S9").Value = 0.0003626415178177816
$ws.Range("T9").Value = 0.0003626415178177817
$ws.Range("I10").Value = 0.04107483513127341
$ws.Range("J10").Value = 0.04107483513127341
$ws.Range("M10").Value = 2.688466333333334
$ws.Range("N10").Value = 8.065399000000001
$ws.Range("O10").Value = 0.1382584398787761
$ws.Range("P10").Value = 0.1382584398787761
$ws.Range("Q10").Value = 0.09593702494955557
$ws.Range("R10").Value = 0.8634332245460001
$ws.Range("S10").Value = 0.005678942623527803
$ws.Range("T10").Value = 0.005678942623527804
$ws.Range("I11").Value = 0.04107483513127341
$ws.Range("J11").Value = 0.04107483513127341
$ws.Range("M11").Value = 0.113399
$ws.Range("N11").Value = 0.340197
$ws.Range("O11").Value = 0.005831714769652435
$ws.Range("P11").Value = 0.005831714769652436
$ws.Range("Q11").Value = 0.004046605515333333
$ws.Range("R11").Value = 0.036419449638
$ws.Range("S11").Value = 0.0002395367226960858
$ws.Range("T11").Value = 0.0002395367226960859
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.011261
$ws.Range("H12").Value = 0.033783
$ws.Range("I12").Value = 0.01296197391260307
$ws.Range("J12").Value = 0.01296197391260308
$ws.Range("M12").Value = 15.47987166666667
$ws.Range("N12").Value = 46.439615
$ws.Range("O12").Value = 0.7960757698994193
$ws.Range("P12").Value = 0.7960757698994194
$ws.Range("Q12").Value = 0.1743188348383334
$ws.Range("R12").Value = 1.568869513545
$ws.Range("S12").Value = 0.01031871336189168
$ws.Range("T12").Value = 0.01031871336189168
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.011261
$ws.Range("H13").Value = 0.033783
$ws.Range("I13").Value = 0.01296197391260307
$ws.Range("J13").Value = 0.01296197391260308
$ws.Range("O13").Value = 0.05100527512565552
$ws.Range("P13").Value = 0.05100527512565553
$ws.Range("Q13").Value = 0.011168761149
$ws.Range("R13").Value = 0.100518850341
$ws.Range("S13").Value = 0.0006611290455838894
$ws.Range("T13").Value = 0.0006611290455838896
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 0.3333333333333333
$ws.Range("G14").Value = 0.011261
$ws.Range("H14").Value = 0.033783
$ws.Range("I14").Value = 0.01296197391260307
$ws.Range("J14").Value = 0.01296197391260308
$ws.Range("M14").Value = 0.171678
$ws.Range("N14").Value = 0.515034
$ws.Range("O14").Value = 0.008828800326496623
$ws.Range("P14").Value = 0.008828800326496624
$ws.Range("Q14").Value = 0.001933265958
$ws.Range("R14").Value = 0.017399393622
$ws.Range("S14").Value = 0.0001144386795116307
$ws.Range("T14").Value = 0.0001144386795116308
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.3333333333333333
$ws.Range("G15").Value = 0.011261
$ws.Range("H15").Value = 0.033783
$ws.Range("I15").Value = 0.01296197391260307
$ws.Range("J15").Value = 0.01296197391260308
$ws.Range("M15").Value = 2.688466333333334
$ws.Range("N15").Value = 8.065399000000001
$ws.Range("O15").Value = 0.1382584398787761
$ws.Range("P15").Value = 0.1382584398787761
$ws.Range("Q15").Value = 0.03027481937966667
$ws.Range("R15").Value = 0.272473374417
$ws.Range("S15").Value = 0.001792102290905896
$ws.Range("T15").Value = 0.001792102290905896
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0.3333333333333333
$ws.Range("G16").Value = 0.011261
$ws.Range("H16").Value = 0.033783
$ws.Range("I16").Value = 0.01296197391260307
$ws.Range("J16").Value = 0.01296197391260308
$ws.Range("M16").Value = 0.113399
$ws.Range("N16").Value = 0.340197
$ws.Range("O16").Value = 0.005831714769652435
$ws.Range("P16").Value = 0.005831714769652436
$ws.Range("Q16").Value = 0.001276986139
$ws.Range("R16").Value = 0.011492875251
$ws.Range("S16").Value = 0.00007559053470997691
$ws.Range("T16").Value = 0.00007559053470997692
$ws.Range("G17").Value = 0.7979533333333334
$ws.Range("H17").Value = 2.39386
$ws.Range("I17").Value = 0.9184841745974011
$ws.Range("J17").Value = 0.9184841745974011
$ws.Range("M17").Value = 15.47987166666667
$ws.Range("N17").Value = 46.439615
$ws.Range("O17").Value = 0.7960757698994193
$ws.Range("P17").Value = 0.7960757698994194
$ws.Range("Q17").Value = 12.35221519598889
$ws.Range("R17").Value = 111.1699367639
$ws.Range("S17").Value = 0.7311829964330587
$ws.Range("T17").Value = 0.7311829964330588
$ws.Range("G18").Value = 0.7979533333333334
$ws.Range("H18").Value = 2.39386
$ws.Range("I18").Value = 0.9184841745974011
$ws.Range("J18").Value = 0.9184841745974011
$ws.Range("O18").Value = 0.05100527512565552
$ws.Range("P18").Value = 0.05100527512565553
$ws.Range("Q18").Value = 0.79141729758
$ws.Range("R18").Value = 7.12275567822
$ws.Range("S18").Value = 0.04684753802390106
$ws.Range("T18").Value = 0.04684753802390107
$ws.Range("G19").Value = 0.7979533333333334
$ws.Range("H19").Value = 2.39386
$ws.Range("I19").Value = 0.9184841745974011
$ws.Range("J19").Value = 0.9184841745974011
$ws.Range("M19").Value = 0.171678
$ws.Range("N19").Value = 0.515034
$ws.Range("O19").Value = 0.008828800326496623
$ws.Range("P19").Value = 0.008828800326496624
$ws.Range("Q19").Value = 0.13699103236
$ws.Range("R19").Value = 1.23291929124
$ws.Range("S19").Value = 0.008109113380567516
$ws.Range("T19").Value = 0.008109113380567518
$ws.Range("G20").Value = 0.7979533333333334
$ws.Range("H20").Value = 2.39386
$ws.Range("I20").Value = 0.9184841745974011
$ws.Range("J20").Value = 0.9184841745974011
$ws.Range("M20").Value = 2.688466333333334
$ws.Range("N20").Value = 8.065399000000001
$ws.Range("O20").Value = 0.1382584398787761
$ws.Range("P20").Value = 0.1382584398787761
$ws.Range("Q20").Value = 2.145270672237778
$ws.Range("R20").Value = 19.30743605014001
$ws.Range("S20").Value = 0.126988189033182
$ws.Range("T20").Value = 0.126988189033182
$ws.Range("G21").Value = 0.7979533333333334
$ws.Range("H21").Value = 2.39386
$ws.Range("I21").Value = 0.9184841745974011
$ws.Range("J21").Value = 0.9184841745974011
$ws.Range("M21").Value = 0.113399
$ws.Range("N21").Value = 0.340197
$ws.Range("O21").Value = 0.005831714769652435
$ws.Range("P21").Value = 0.005831714769652436
$ws.Range("Q21").Value = 0.09048711004666668
$ws.Range("R21").Value = 0.8143839904200001
$ws.Range("S21").Value = 0.00535633772669169
$ws.Range("T21").Value = 0.00535633772669169
